$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the "Play 3 Lucky Leprechauns Free: Game Review" heading (it is
#    paragraph #2 in the document).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph reading "Play 3 Lucky Leprechauns Free:
#    Game Review" right before the trailing "Prompt: ..." paragraph (i.e.
#    right after what is now the last "What we don't like" bullet point).
$count = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($count - 1)
$prevPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($count)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 3 Lucky Leprechauns Free: Game Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml)

# 3. Replace the old "Prompt: ..." text (now the last paragraph) with the
#    meta-description sentence, keeping its italic run formatting intact.
$oldText = 'Prompt: Create a feature image fitting the game "3 Lucky Leprechauns". The image should be in cartoon style and feature a happy Maya warrior with glasses. Sorry, but the prompt does not fit the context of the given game. "3 Lucky Leprechauns" is an Irish-themed slot game, and a happy Maya warrior with glasses is not relevant to the game''s theme. Please provide a suitable prompt that would aptly capture the essence of the game.'
$newText = 'Read our 3 Lucky Leprechauns game review and play for free. Win big with bonus games and Irish folklore theme.'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
